# Corrijo datos y graficas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (XGB) - corrected values
$ws.Range("B6").Value = 442.0
$ws.Range("C6").Value = 67.0
$ws.Range("D6").Value = 343.0
$ws.Range("E6").Value = 66.0
$ws.Range("F6").Value = 0.8683693516699411
$ws.Range("G6").Value = 0.8700787401574803
$ws.Range("H6").Value = 0.8365853658536585
$ws.Range("I6").Value = 0.8692232055063914
$ws.Range("J6").Value = 0.855119825708061
$ws.Range("K6").Value = 0.7068303950900682
$ws.Range("L6").Value = 0.8531677098649104
$ws.Range("M6").Value = 0.9190896869598638

# Row 7 (k-NN) - corrected values
$ws.Range("C7").Value = 83.0
$ws.Range("D7").Value = 327.0
$ws.Range("F7").Value = 0.8442776735459663
$ws.Range("H7").Value = 0.7975609756097561
$ws.Range("I7").Value = 0.8645533141210374
$ws.Range("J7").Value = 0.8464052287581699
$ws.Range("K7").Value = 0.6874308758542416
$ws.Range("L7").Value = 0.8405360576568032
$ws.Range("M7").Value = 0.8416938736316497
